$wb = $excel.ActiveWorkbook

# Sheet ALC, row 15
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1977.6578
$ws.Range("I15").Value = 1977.6578
$ws.Range("K15").Value = 5932.9734
$ws.Range("M15").Value = -5763.9734

# Sheet ALC, row 29
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 74.5
$ws.Range("I29").Value = 74.5
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 223.5
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = 57.5
$ws.Range("N29").Value = $null

# Sheet ALC, row 41
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 259.5
$ws.Range("I41").Value = 280.66666
$ws.Range("J41").Value = 238.33333
$ws.Range("K41").Value = 280.66666
$ws.Range("L41").Value = 238.33333
$ws.Range("M41").Value = 159.33334
$ws.Range("N41").Value = -1118.33333

# Sheet ALC, row 103
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 893.8
$ws.Range("J103").Value = 893.8
$ws.Range("L103").Value = 2681.4
$ws.Range("N103").Value = -3853.4

# Sheet ALC, row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1403.4849
$ws.Range("I112").Value = 593.2222
$ws.Range("J112").Value = 1707.3334
$ws.Range("K112").Value = 1779.6666
$ws.Range("L112").Value = 5122.0002
$ws.Range("M112").Value = -671.6666
$ws.Range("N112").Value = -7338.0002

# Sheet ALC, row 131
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 4423.857
$ws.Range("I131").Value = 4423.857
$ws.Range("K131").Value = 13271.571
$ws.Range("M131").Value = -8231.571

# Sheet ALC, row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 917.0833
$ws.Range("I132").Value = 924.65
$ws.Range("J132").Value = 879.25
$ws.Range("K132").Value = 2773.95
$ws.Range("L132").Value = 2637.75
$ws.Range("M132").Value = -243.9499999999998
$ws.Range("N132").Value = -7697.75

# Sheet ALC, row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3520.2144
$ws.Range("I137").Value = 1910.5625
$ws.Range("K137").Value = 5731.6875
$ws.Range("M137").Value = -3181.6875

# Sheet ALC, row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4692.619
$ws.Range("I138").Value = 1601.8462
$ws.Range("J138").Value = 9715.125
$ws.Range("K138").Value = 4805.5386
$ws.Range("L138").Value = 29145.375
$ws.Range("M138").Value = 334.4614000000001
$ws.Range("N138").Value = -39425.375

# Sheet ARM, row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 986.8889
$ws.Range("I2").Value = 697.4286
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = 697.4286
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = -584.4286
$ws.Range("N2").Value = -2226

# Sheet ARM, row 10
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").Value = $null

# Sheet ARM, row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3232.5527
$ws.Range("I32").Value = 2611
$ws.Range("J32").Value = 8515.75
$ws.Range("K32").Value = 2611
$ws.Range("L32").Value = 8515.75
$ws.Range("M32").Value = -2324
$ws.Range("N32").Value = -9089.75

# Sheet ARM, row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2996
$ws.Range("I61").Value = 2995.5
$ws.Range("K61").Value = 2995.5
$ws.Range("M61").Value = -2783.5

# Sheet ARM, row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2201.0454
$ws.Range("I74").Value = 2382.1052
$ws.Range("J74").Value = 1054.3334
$ws.Range("K74").Value = 2382.1052
$ws.Range("L74").Value = 1054.3334
$ws.Range("M74").Value = -1508.1052
$ws.Range("N74").Value = -2802.3334

# Sheet ARM, row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2201.0454
$ws.Range("I77").Value = 2382.1052
$ws.Range("J77").Value = 1054.3334
$ws.Range("K77").Value = 11910.526
$ws.Range("L77").Value = 5271.666999999999
$ws.Range("M77").Value = -7542.526
$ws.Range("N77").Value = -14007.667

# Sheet ARM, row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1566.8889
$ws.Range("I102").Value = 1566.8889
$ws.Range("K102").Value = 1566.8889
$ws.Range("M102").Value = 55.11110000000008

# Sheet ARM, row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1139.6428
$ws.Range("I110").Value = 1063.9
$ws.Range("K110").Value = 1063.9
$ws.Range("M110").Value = 981.0999999999999

# Sheet ARM, row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 986.8889
$ws.Range("I116").Value = 697.4286
$ws.Range("J116").Value = 2000
$ws.Range("K116").Value = 697.4286
$ws.Range("L116").Value = 2000
$ws.Range("M116").Value = 1596.5714
$ws.Range("N116").Value = -6588

# Sheet ARM, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1654.9546
$ws.Range("I132").Value = 1281.6757
$ws.Range("K132").Value = 3845.0271
$ws.Range("M132").Value = -1315.0271

# Sheet ARM, row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2996
$ws.Range("I136").Value = 2995.5
$ws.Range("K136").Value = 8986.5
$ws.Range("M136").Value = -6436.5

# Sheet BSM, row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 986.8889
$ws.Range("I3").Value = 697.4286
$ws.Range("J3").Value = 2000
$ws.Range("K3").Value = 697.4286
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = -583.4286
$ws.Range("N3").Value = -2228

# Sheet BSM, row 22
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 643.25
$ws.Range("J22").Value = 595
$ws.Range("L22").Value = 595
$ws.Range("N22").Value = -941

# Sheet BSM, row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 269.0909
$ws.Range("I94").Value = 307.75
$ws.Range("J94").Value = 166
$ws.Range("K94").Value = 307.75
$ws.Range("L94").Value = 166
$ws.Range("M94").Value = 143.25
$ws.Range("N94").Value = -1068

# Sheet BSM, row 135
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 47949
$ws.Range("J135").Value = 47949
$ws.Range("L135").Value = 47949
$ws.Range("N135").Value = -58089

# Sheet CRP, row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3279.7222
$ws.Range("I31").Value = 2281.5
$ws.Range("J31").Value = 4527.5
$ws.Range("K31").Value = 2281.5
$ws.Range("L31").Value = 4527.5
$ws.Range("M31").Value = -1986.5
$ws.Range("N31").Value = -5117.5

# Sheet CRP, row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3279.7222
$ws.Range("I34").Value = 2281.5
$ws.Range("J34").Value = 4527.5
$ws.Range("K34").Value = 2281.5
$ws.Range("L34").Value = 4527.5
$ws.Range("M34").Value = -2079.5
$ws.Range("N34").Value = -4931.5

# Sheet CRP, row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2096.2856
$ws.Range("I58").Value = 1871.1111
$ws.Range("K58").Value = 1871.1111
$ws.Range("M58").Value = -1668.1111

# Sheet CRP, row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 5000
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").Value = $null

# Sheet CRP, row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2096.2856
$ws.Range("I136").Value = 1871.1111
$ws.Range("K136").Value = 5613.3333
$ws.Range("M136").Value = -3063.3333

# Sheet GSM, row 43
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 26252.5
$ws.Range("J43").Value = 26252.5
$ws.Range("L43").Value = 26252.5
$ws.Range("N43").Value = -26554.5

# Sheet GSM, row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4613.857
$ws.Range("I102").Value = 4613.857
$ws.Range("K102").Value = 4613.857
$ws.Range("M102").Value = -2991.857

# Sheet GSM, row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 660.5
$ws.Range("I122").Value = 638.5333000000001
$ws.Range("J122").Value = 770.3333
$ws.Range("K122").Value = 1915.5999
$ws.Range("L122").Value = 2310.9999
$ws.Range("M122").Value = 534.4000999999998
$ws.Range("N122").Value = -7210.9999

# Sheet GSM, row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2697.6
$ws.Range("I132").Value = 2343.5386
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 7030.6158
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -4500.6158
$ws.Range("N132").Value = -20057

# Sheet LTW, row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6403.9
$ws.Range("I40").Value = 6403.9
$ws.Range("K40").Value = 6403.9
$ws.Range("M40").Value = -6267.9

# Sheet LTW, row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3284.7144
$ws.Range("I68").Value = 3284.7144
$ws.Range("K68").Value = 3284.7144
$ws.Range("M68").Value = -2535.7144

# Sheet LTW, row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 3284.7144
$ws.Range("I71").Value = 3284.7144
$ws.Range("K71").Value = 16423.572
$ws.Range("M71").Value = -12679.572

# Sheet LTW, row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2000
$ws.Range("I82").Value = 2000
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 2000
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -1639
$ws.Range("N82").Value = $null

# Sheet LTW, row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2000
$ws.Range("I85").Value = 2000
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 2000
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -752
$ws.Range("N85").Value = $null

# Sheet LTW, row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4022.1765
$ws.Range("I136").Value = 4023.5625
$ws.Range("K136").Value = 12070.6875
$ws.Range("M136").Value = -9520.6875

# Sheet WVR, row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1408.9
$ws.Range("I122").Value = 1503.7778
$ws.Range("J122").Value = 555
$ws.Range("K122").Value = 4511.3334
$ws.Range("L122").Value = 1665
$ws.Range("M122").Value = -2061.3334
$ws.Range("N122").Value = -6565

# Sheet WVR, row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2586.2856
$ws.Range("J126").Value = 2340
$ws.Range("L126").Value = 7020
$ws.Range("N126").Value = -11960

# Sheet WVR, row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2661.5715
$ws.Range("I132").Value = 1871.4286
$ws.Range("J132").Value = 4241.857
$ws.Range("K132").Value = 5614.2858
$ws.Range("L132").Value = 12725.571
$ws.Range("M132").Value = -3084.2858
$ws.Range("N132").Value = -17785.571

# Sheet WVR, row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 827.2308
$ws.Range("I136").Value = 832.4545000000001
$ws.Range("J136").Value = 798.5
$ws.Range("K136").Value = 2497.3635
$ws.Range("L136").Value = 2395.5
$ws.Range("M136").Value = 52.63649999999961
$ws.Range("N136").Value = -7495.5
